$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking strings
# (e.g. "1.000", "0.000007939") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.840.65'
$ws.Range("E2").Value = '  -0.27%  '

# Row 3
$ws.Range("D3").Value = '1.875.71'
$ws.Range("E3").Value = '  +0.00%  '

# Row 4
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").Value = '0.7215'
$ws.Range("E5").Value = '  -2.46%  '

# Row 6
$ws.Range("D6").Value = '241.93'
$ws.Range("E6").Value = '  -0.29%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '0.3147'
$ws.Range("E8").Value = '  -0.52%  '

# Row 9
$ws.Range("D9").Value = '0.07582'
$ws.Range("E9").Value = '  +5.15%  '

# Row 10
$ws.Range("D10").Value = '24.57'
$ws.Range("E10").Value = '  -0.60%  '

# Row 11
$ws.Range("D11").Value = '0.08187'
$ws.Range("E11").Value = '  -2.70%  '

# Row 12
$ws.Range("D12").Value = '1.912.19'
$ws.Range("E12").Value = '  +1.57%  '

# Row 13
$ws.Range("D13").Value = '0.7454'
$ws.Range("E13").Value = '  -0.83%  '

# Row 14
$ws.Range("D14").Value = '5.330'
$ws.Range("E14").Value = '  -1.81%  '

# Row 15
$ws.Range("D15").Value = '92.53'
$ws.Range("E15").Value = '  -0.08%  '

# Row 16
$ws.Range("D16").Value = '30.103.42'
$ws.Range("E16").Value = '  +0.57%  '

# Row 17
$ws.Range("D17").Value = '6.013'
$ws.Range("E17").Value = '  -1.53%  '

# Row 18
$ws.Range("D18").Value = '247.22'
$ws.Range("E18").Value = '  +1.51%  '

# Row 19
$ws.Range("D19").Value = '0.000007939'
$ws.Range("E19").Value = '  +1.52%  '

# Row 20
$ws.Range("D20").Value = '13.44'
$ws.Range("E20").Value = '  -1.05%  '

# Row 21
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.193.39'
$ws.Range("E21").Value = '  +2.00%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.004'
$ws.Range("E22").Value = '  +0.43%  '

# Row 23
$ws.Range("E23").Value = '  +0.13%  '

# Row 24
$ws.Range("D24").Value = '7.747'
$ws.Range("E24").Value = '  -3.25%  '

# Row 25
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.274'
$ws.Range("E25").Value = '  +0.02%  '

# Row 26
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1502'
$ws.Range("E26").Value = '  -3.36%  '

# Row 27
$ws.Range("D27").Value = '164.35'
$ws.Range("E27").Value = '  -0.41%  '

# Row 28
$ws.Range("D28").Value = '18.61'
$ws.Range("E28").Value = '  -0.01%  '

# Row 29
$ws.Range("D29").Value = '2.009'
$ws.Range("E29").Value = '  -1.32%  '

# Row 30
$ws.Range("D30").Value = '1.437'
$ws.Range("E30").Value = '  -3.89%  '

# Row 31
$ws.Range("D31").Value = '4.535'
$ws.Range("E31").Value = '  -1.43%  '

# Row 32
$ws.Range("D32").Value = '1.525'
$ws.Range("E32").Value = '  -0.27%  '

# Row 33
$ws.Range("D33").Value = '4.204'
$ws.Range("E33").Value = '  -1.10%  '

# Row 34
$ws.Range("D34").Value = '0.05438'
$ws.Range("E34").Value = '  +2.29%  '

# Row 35
$ws.Range("D35").Value = '1.230'
$ws.Range("E35").Value = '  -0.59%  '

# Row 36
$ws.Range("D36").Value = '0.7412'
$ws.Range("E36").Value = '  -1.74%  '

# Row 37
$ws.Range("D37").Value = '1.006'
$ws.Range("E37").Value = '  +0.31%  '

# Row 38
$ws.Range("D38").Value = '2.697'
$ws.Range("E38").Value = '  +0.05%  '

# Row 39
$ws.Range("D39").Value = '0.01928'
$ws.Range("E39").Value = '  -1.14%  '

# Row 40
$ws.Range("D40").Value = '2.742'
$ws.Range("E40").Value = '  -0.35%  '

# Row 41
$ws.Range("D41").Value = '0.4464'
$ws.Range("E41").Value = '  -1.37%  '

# Row 42
$ws.Range("D42").Value = '0.8836'
$ws.Range("E42").Value = '  +3.34%  '

# Row 43
$ws.Range("D43").Value = '5.999'
$ws.Range("E43").Value = '  -1.15%  '

# Row 44
$ws.Range("D44").Value = '71.84'
$ws.Range("E44").Value = '  -0.69%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.040.99'
$ws.Range("E45").Value = '  -6.48%  '

# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '104.05'
$ws.Range("E46").Value = '  +0.45%  '

# Row 47
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.02%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.477'
$ws.Range("E48").Value = '  -2.26%  '

# Row 49
$ws.Range("D49").Value = '1.817'
$ws.Range("E49").Value = '  -1.18%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.643'
$ws.Range("E50").Value = '  +1.09%  '

# Row 51
$ws.Range("D51").Value = '2.055.59'
$ws.Range("E51").Value = '  +1.72%  '

